$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as text (matches source data which uses
# localized/grouped numeric strings such as "28.138.71" and plain decimals like
# "311.42"). Without forcing text format, Excel auto-converts plain decimals to
# numbers, which would not match the original text values.
$priceCells = @("D2","D3","D5","D7","D8","D9","D10","D11","D12","D13","D14","D16","D17","D18","D19","D20","D22","D23","D24","D25","D26","D27","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D49","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.138.71'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '1.866.55'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '311.42'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.5039'
$ws.Range("E7").Value = '  -1.56%  '
$ws.Range("D8").Value = '0.3917'
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").Value = '0.09676'
$ws.Range("E9").Value = '  -4.26%  '
$ws.Range("D10").Value = '1.137'
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("D11").Value = '40.92'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = '6.494'
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '20.89'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '1.875.08'
$ws.Range("E14").Value = '  +2.29%  '
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").Value = '7.394'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '0.00001126'
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("D18").Value = '92.83'
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("D19").Value = '0.06609'
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").Value = '17.51'
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '6.158'
$ws.Range("E22").Value = '  +1.56%  '
$ws.Range("D23").Value = '28.218.18'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '11.32'
$ws.Range("E24").Value = '  +1.12%  '
$ws.Range("D25").Value = '2.281'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("D26").Value = '2.529'
$ws.Range("E26").Value = '  +2.21%  '
$ws.Range("D27").Value = '2.082.11'
$ws.Range("E27").Value = '  +2.36%  '
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("D29").Value = '158.14'
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").Value = '127.36'
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("E31").Value = '  -3.34%  '
$ws.Range("D32").Value = '1.063'
$ws.Range("E32").Value = '  -0.84%  '
$ws.Range("D33").Value = '5.627'
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("D34").Value = '3.618'
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("D35").Value = '9.546'
$ws.Range("E35").Value = '  +4.48%  '
$ws.Range("D36").Value = '0.06717'
$ws.Range("E36").Value = '  -2.87%  '
$ws.Range("D37").Value = '0.02386'
$ws.Range("E37").Value = '  +1.58%  '
$ws.Range("D38").Value = '0.2174'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = '11.48'
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("D40").Value = '0.6348'
$ws.Range("E40").Value = '  +1.26%  '
$ws.Range("D41").Value = '4.966'
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("D42").Value = '1.177'
$ws.Range("E42").Value = '  +1.32%  '
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '13.58'
$ws.Range("E44").Value = '  +1.80%  '
$ws.Range("D45").Value = '0.6001'
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("E46").Value = '  -1.59%  '
$ws.Range("D47").Value = '1.258'
$ws.Range("E47").Value = '  -2.37%  '
$ws.Range("D48").Value = '124.14'
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("D49").Value = '1.987'
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").Value = '0.06831'
$ws.Range("E51").Value = '  +0.70%  '
